# TC11 - Create and Manage API Keys: prepend a new "button_apiKeyManagement_
# trNthChild" column ahead of the existing "input_KeyName" column, and record
# a row index ("2") under the new header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing column (header "input_KeyName" + its data) from A to B
# by inserting a new blank column at A.
$ws.Columns.Item(1).Insert()

# New column A header.
$ws.Range("A1").Value = "button_apiKeyManagement_trNthChild"

# Give A1 the same look as B1 (bold / bordered / centered "Pandas" header
# style) by copying formats over instead of assigning .Style directly, so we
# reuse the existing style record rather than minting a near-duplicate one.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)   # xlPasteFormats

# New column A data value "2", written as text (matching the original
# inlineStr cell type) rather than a number. Stage it in a scratch cell
# formatted as Text, then copy only the *value* into A2 so A2's own style
# stays at the sheet default (no stray number-format style is created).
$ws.Range("D1").NumberFormat = "@"
$ws.Range("D1").Value = "2"
$ws.Range("D1").Copy()
$ws.Range("A2").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("D1").Clear()

# Column widths: A=36, B=15. The engine's ColumnWidth setter stores a value
# 5/6 wider than requested (fixed padding offset), so compensate to land on
# the exact target widths.
$ws.Columns.Item(1).ColumnWidth = 36 - 5/6
$ws.Columns.Item(2).ColumnWidth = 15 - 5/6
